# Add a "last updated" style date stamp to the About sheet (cell C1),
# matching the author's commit that added a date value next to the title.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("About")

$cell = $ws.Range("C1")
$cell.Value = 44307
$cell.NumberFormat = "mm-dd-yy"
